# Auto-generated edit script: updates market-price-derived columns (H-N)
# on each job sheet, per the scheduled-runner refresh described in the commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 95.59999999999999
$ws.Range("I9").Value = 95.59999999999999
$ws.Range("K9").Value = 95.59999999999999
$ws.Range("M9").Value = 73.40000000000001

$ws.Range("H15").Value = 626260.1
$ws.Range("I15").Value = 626260.1
$ws.Range("K15").Value = 1878780.3
$ws.Range("M15").Value = -1878611.3

$ws.Range("H80").Value = 59685.543
$ws.Range("I80").Value = 119172.82
$ws.Range("J80").Value = 3503.111
$ws.Range("K80").Value = 357518.46
$ws.Range("L80").Value = 10509.333
$ws.Range("M80").Value = -356520.46
$ws.Range("N80").Value = -12505.333

$ws.Range("H83").Value = 59685.543
$ws.Range("I83").Value = 119172.82
$ws.Range("J83").Value = 3503.111
$ws.Range("K83").Value = 1072555.38
$ws.Range("L83").Value = 31527.999
$ws.Range("M83").Value = -1067563.38
$ws.Range("N83").Value = -41511.999

$ws.Range("H98").Value = 63798.89
$ws.Range("I98").Value = 63798.89
$ws.Range("K98").Value = 63798.89
$ws.Range("M98").Value = -62300.89

$ws.Range("H111").Value = 1184
$ws.Range("I111").Value = 1020.5714
$ws.Range("K111").Value = 3061.7142
$ws.Range("M111").Value = 5.285799999999654

$ws.Range("H122").Value = 63798.89
$ws.Range("I122").Value = 63798.89
$ws.Range("K122").Value = 191396.67
$ws.Range("M122").Value = -188946.67

$ws.Range("H137").Value = 11131.064
$ws.Range("I137").Value = 16724.945
$ws.Range("K137").Value = 50174.835
$ws.Range("M137").Value = -47624.835

$ws.Range("H138").Value = 196954.2
$ws.Range("I138").Value = 878776.25
$ws.Range("J138").Value = 3985.698
$ws.Range("K138").Value = 2636328.75
$ws.Range("L138").Value = 11957.094
$ws.Range("M138").Value = -2631188.75
$ws.Range("N138").Value = -22237.094

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11308.269
$ws.Range("I32").Value = 11308.269
$ws.Range("K32").Value = 11308.269
$ws.Range("M32").Value = -11021.269

$ws.Range("H46").Value = 8778.666999999999
$ws.Range("J46").Value = 8001
$ws.Range("L46").Value = 8001
$ws.Range("N46").Value = -8639

$ws.Range("H61").Value = 10330.036
$ws.Range("I61").Value = 13410.25
$ws.Range("K61").Value = 13410.25
$ws.Range("M61").Value = -13198.25

$ws.Range("H74").Value = 1631.1428
$ws.Range("I74").Value = 615.5
$ws.Range("K74").Value = 615.5
$ws.Range("M74").Value = 258.5

$ws.Range("H77").Value = 1631.1428
$ws.Range("I77").Value = 615.5
$ws.Range("K77").Value = 3077.5
$ws.Range("M77").Value = 1290.5

$ws.Range("H88").Value = 1704.0667
$ws.Range("J88").Value = 1716.8334
$ws.Range("L88").Value = 1716.8334
$ws.Range("N88").Value = -2528.8334

$ws.Range("H91").Value = 1704.0667
$ws.Range("J91").Value = 1716.8334
$ws.Range("L91").Value = 1716.8334
$ws.Range("N91").Value = -4524.8334

$ws.Range("H122").Value = 861751.4399999999
$ws.Range("I122").Value = 4538.5
$ws.Range("K122").Value = 13615.5
$ws.Range("M122").Value = -11165.5

$ws.Range("H136").Value = 10330.036
$ws.Range("I136").Value = 13410.25
$ws.Range("K136").Value = 40230.75
$ws.Range("M136").Value = -37680.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 285.0476
$ws.Range("J80").Value = 234.2
$ws.Range("L80").Value = 234.2
$ws.Range("N80").Value = -2230.2

$ws.Range("H83").Value = 285.0476
$ws.Range("J83").Value = 234.2
$ws.Range("L83").Value = 1171
$ws.Range("N83").Value = -11155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 219.6
$ws.Range("I7").Value = 212.33333
$ws.Range("J7").Value = 230.5
$ws.Range("K7").Value = 212.33333
$ws.Range("L7").Value = 230.5
$ws.Range("M7").Value = -99.33332999999999
$ws.Range("N7").Value = -456.5

$ws.Range("H107").Value = 83343910
$ws.Range("I107").Value = 125015000
$ws.Range("K107").Value = 125015000
$ws.Range("M107").Value = -125013080

$ws.Range("H132").Value = 1826.4445
$ws.Range("I132").Value = 1572.6
$ws.Range("K132").Value = 4717.799999999999
$ws.Range("M132").Value = -2187.799999999999

$ws.Range("H141").Value = 171408.19
$ws.Range("J141").Value = 183863.75
$ws.Range("L141").Value = 183863.75
$ws.Range("N141").Value = -194223.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3688.157
$ws.Range("I122").Value = 967.75
$ws.Range("J122").Value = 4931.7715
$ws.Range("K122").Value = 8709.75
$ws.Range("L122").Value = 44385.9435
$ws.Range("M122").Value = -6259.75
$ws.Range("N122").Value = -49285.9435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7640.476
$ws.Range("I102").Value = 8141.9443
$ws.Range("K102").Value = 8141.9443
$ws.Range("M102").Value = -6519.9443

$ws.Range("H122").Value = 10430.2
$ws.Range("I122").Value = 7671.7144
$ws.Range("J122").Value = 16866.666
$ws.Range("K122").Value = 23015.1432
$ws.Range("L122").Value = 50599.99800000001
$ws.Range("M122").Value = -20565.1432
$ws.Range("N122").Value = -55499.99800000001

$ws.Range("H132").Value = 2097.55
$ws.Range("I132").Value = 2030.6875
$ws.Range("J132").Value = 2365
$ws.Range("K132").Value = 6092.0625
$ws.Range("L132").Value = 7095
$ws.Range("M132").Value = -3562.0625
$ws.Range("N132").Value = -12155

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1419.375
$ws.Range("I16").Value = 1265.0938
$ws.Range("K16").Value = 1265.0938
$ws.Range("M16").Value = -1095.0938

$ws.Range("H22").Value = 1250
$ws.Range("J22").Value = 1666.6666
$ws.Range("L22").Value = 1666.6666
$ws.Range("N22").Value = -2256.6666

$ws.Range("H27").Value = 1250
$ws.Range("J27").Value = 1666.6666
$ws.Range("L27").Value = 1666.6666
$ws.Range("N27").Value = -1880.6666

$ws.Range("H40").Value = 19289.883
$ws.Range("I40").Value = 25998.15
$ws.Range("J40").Value = 9706.643
$ws.Range("K40").Value = 25998.15
$ws.Range("L40").Value = 9706.643
$ws.Range("M40").Value = -25862.15
$ws.Range("N40").Value = -9978.643

$ws.Range("H122").Value = 5049.6665
$ws.Range("I122").Value = 4303.1035
$ws.Range("K122").Value = 12909.3105
$ws.Range("M122").Value = -10459.3105

$ws.Range("H132").Value = 786561
$ws.Range("I132").Value = 994794
$ws.Range("J132").Value = 5687.25
$ws.Range("K132").Value = 2984382
$ws.Range("L132").Value = 17061.75
$ws.Range("M132").Value = -2981852
$ws.Range("N132").Value = -22121.75

$ws.Range("H136").Value = 5155.4053
$ws.Range("I136").Value = 3852.5
$ws.Range("K136").Value = 11557.5
$ws.Range("M136").Value = -9007.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 100000
$ws.Range("J120").Value = 100000
$ws.Range("L120").Value = 100000
$ws.Range("N120").Value = -109676

$ws.Range("H122").Value = 13289.54
$ws.Range("I122").Value = 2138.7632
$ws.Range("J122").Value = 48600.332
$ws.Range("K122").Value = 6416.2896
$ws.Range("L122").Value = 145800.996
$ws.Range("M122").Value = -3966.2896
$ws.Range("N122").Value = -150700.996

$ws.Range("H132").Value = 7530.4927
$ws.Range("I132").Value = 8050.654
$ws.Range("K132").Value = 24151.962
$ws.Range("M132").Value = -21621.962

$ws.Range("H136").Value = 499083.6
$ws.Range("I136").Value = 572210.6
$ws.Range("K136").Value = 1716631.8
$ws.Range("M136").Value = -1714081.8
